$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.0421005
$ws.Range("H2").Value = 0.084201
$ws.Range("I2").Value = 0.02355433709362141
$ws.Range("J2").Value = 0.02354456735134313
$ws.Range("M2").Value = 4.296436999999999
$ws.Range("N2").Value = 8.592873999999998
$ws.Range("O2").Value = 0.08737129157293111
$ws.Range("P2").Value = 0.06876644796033347
$ws.Range("Q2").Value = 0.1808821459185
$ws.Range("R2").Value = 0.7235285836739999
$ws.Range("S2").Value = 0.002057972854013902
$ws.Range("T2").Value = 0.001619076265514704

# Row 3
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.0421005
$ws.Range("H3").Value = 0.084201
$ws.Range("I3").Value = 0.02355433709362141
$ws.Range("J3").Value = 0.02354456735134313
$ws.Range("M3").Value = 13.32522833333333
$ws.Range("N3").Value = 39.975685
$ws.Range("O3").Value = 0.2709785829485105
$ws.Range("P3").Value = 0.3199146015909443
$ws.Range("Q3").Value = 0.5609987754475
$ws.Range("R3").Value = 3.365992652685
$ws.Range("S3").Value = 0.006382720887921066
$ws.Range("T3").Value = 0.007532250883836092

# Row 4
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.0421005
$ws.Range("H4").Value = 0.084201
$ws.Range("I4").Value = 0.02355433709362141
$ws.Range("J4").Value = 0.02354456735134313
$ws.Range("M4").Value = 6.89049
$ws.Range("N4").Value = 20.67147
$ws.Range("O4").Value = 0.140123318663899
$ws.Range("P4").Value = 0.1654281868928364
$ws.Range("Q4").Value = 0.290093074245
$ws.Range("R4").Value = 1.74055844547
$ws.Range("S4").Value = 0.003300511882486409
$ws.Range("T4").Value = 0.003894935088108965

# Row 5
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.0421005
$ws.Range("H5").Value = 0.084201
$ws.Range("I5").Value = 0.02355433709362141
$ws.Range("J5").Value = 0.02354456735134313
$ws.Range("M5").Value = 18.2696115
$ws.Range("N5").Value = 36.539223
$ws.Range("O5").Value = 0.3715263492262718
$ws.Range("P5").Value = 0.292413525083752
$ws.Range("Q5").Value = 0.7691597789557499
$ws.Range("R5").Value = 3.076639115823
$ws.Range("S5").Value = 0.008751056868838114
$ws.Range("T5").Value = 0.006884749935778063

# Row 6
$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.0421005
$ws.Range("H6").Value = 0.084201
$ws.Range("I6").Value = 0.02355433709362141
$ws.Range("J6").Value = 0.02354456735134313
$ws.Range("M6").Value = 1.355562
$ws.Range("N6").Value = 4.066686000000001
$ws.Range("O6").Value = 0.02756637715092428
$ws.Range("P6").Value = 0.03254458882907125
$ws.Range("Q6").Value = 0.05706983798100001
$ws.Range("R6").Value = 0.342419027886
$ws.Range("S6").Value = 0.0006493077398627735
$ws.Range("T6").Value = 0.0007662482636078371

# Row 7
$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.0421005
$ws.Range("H7").Value = 0.084201
$ws.Range("I7").Value = 0.02355433709362141
$ws.Range("J7").Value = 0.02354456735134313
$ws.Range("M7").Value = 5.037141666666667
$ws.Range("N7").Value = 15.111425
$ws.Range("O7").Value = 0.1024340804374633
$ws.Range("P7").Value = 0.1209326496430627
$ws.Range("Q7").Value = 0.2120661827375
$ws.Range("R7").Value = 1.272397096425
$ws.Range("S7").Value = 0.002412766860499141
$ws.Range("T7").Value = 0.00284730691449747

# Row 8
$ws.Range("E8").Value = 1.0
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.001483333333333333
$ws.Range("H8").Value = 0.00445
$ws.Range("I8").Value = 0.000829893548901757
$ws.Range("J8").Value = 0.001244323995124487
$ws.Range("M8").Value = 4.296436999999999
$ws.Range("N8").Value = 8.592873999999998
$ws.Range("O8").Value = 0.08737129157293111
$ws.Range("P8").Value = 0.06876644796033347
$ws.Range("Q8").Value = 0.006373048216666665
$ws.Range("R8").Value = 0.03823828929999999
$ws.Range("S8").Value = 0.00007250887123558997
$ws.Range("T8").Value = 0.00008556774125652227

# Row 9
$ws.Range("E9").Value = 1.0
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.001483333333333333
$ws.Range("H9").Value = 0.00445
$ws.Range("I9").Value = 0.000829893548901757
$ws.Range("J9").Value = 0.001244323995124487
$ws.Range("M9").Value = 13.32522833333333
$ws.Range("N9").Value = 39.975685
$ws.Range("O9").Value = 0.2709785829485105
$ws.Range("P9").Value = 0.3199146015909443
$ws.Range("Q9").Value = 0.01976575536111111
$ws.Range("R9").Value = 0.17789179825
$ws.Range("S9").Value = 0.0002248833778795085
$ws.Range("T9").Value = 0.0003980774151503023

# Row 10
$ws.Range("E10").Value = 1.0
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.001483333333333333
$ws.Range("H10").Value = 0.00445
$ws.Range("I10").Value = 0.000829893548901757
$ws.Range("J10").Value = 0.001244323995124487
$ws.Range("M10").Value = 6.89049
$ws.Range("N10").Value = 20.67147
$ws.Range("O10").Value = 0.140123318663899
$ws.Range("P10").Value = 0.1654281868928364
$ws.Range("Q10").Value = 0.0102208935
$ws.Range("R10").Value = 0.09198804149999999
$ws.Range("S10").Value = 0.0001162874382098749
$ws.Range("T10").Value = 0.0002058462624206945

# Row 11
$ws.Range("E11").Value = 1.0
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.001483333333333333
$ws.Range("H11").Value = 0.00445
$ws.Range("I11").Value = 0.000829893548901757
$ws.Range("J11").Value = 0.001244323995124487
$ws.Range("M11").Value = 18.2696115
$ws.Range("N11").Value = 36.539223
$ws.Range("O11").Value = 0.3715263492262718
$ws.Range("P11").Value = 0.292413525083752
$ws.Range("Q11").Value = 0.027099923725
$ws.Range("R11").Value = 0.16259954235
$ws.Range("S11").Value = 0.0003083273204699043
$ws.Range("T11").Value = 0.0003638571657606486

# Row 12
$ws.Range("E12").Value = 1.0
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.001483333333333333
$ws.Range("H12").Value = 0.00445
$ws.Range("I12").Value = 0.000829893548901757
$ws.Range("J12").Value = 0.001244323995124487
$ws.Range("M12").Value = 1.355562
$ws.Range("N12").Value = 4.066686000000001
$ws.Range("O12").Value = 0.02756637715092428
$ws.Range("P12").Value = 0.03254458882907125
$ws.Range("Q12").Value = 0.0020107503
$ws.Range("R12").Value = 0.0180967527
$ws.Range("S12").Value = 0.00002287715856414486
$ws.Range("T12").Value = 0.00004049601279147368

# Row 13
$ws.Range("E13").Value = 1.0
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.001483333333333333
$ws.Range("H13").Value = 0.00445
$ws.Range("I13").Value = 0.000829893548901757
$ws.Range("J13").Value = 0.001244323995124487
$ws.Range("M13").Value = 5.037141666666667
$ws.Range("N13").Value = 15.111425
$ws.Range("O13").Value = 0.1024340804374633
$ws.Range("P13").Value = 0.1209326496430627
$ws.Range("Q13").Value = 0.007471760138888889
$ws.Range("R13").Value = 0.06724584125000001
$ws.Range("S13").Value = 0.00008500938254273447
$ws.Range("T13").Value = 0.0001504793977448456

# Row 14
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 1.743794
$ws.Range("H14").Value = 3.487588
$ws.Range("I14").Value = 0.9756157693574768
$ws.Range("J14").Value = 0.9752111086535323
$ws.Range("M14").Value = 4.296436999999999
$ws.Range("N14").Value = 8.592873999999998
$ws.Range("O14").Value = 0.08737129157293111
$ws.Range("P14").Value = 0.06876644796033347
$ws.Range("Q14").Value = 7.492101061977999
$ws.Range("R14").Value = 29.96840424791199
$ws.Range("S14").Value = 0.08524080984768162
$ws.Range("T14").Value = 0.06706180395356225

# Row 15
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 1.743794
$ws.Range("H15").Value = 3.487588
$ws.Range("I15").Value = 0.9756157693574768
$ws.Range("J15").Value = 0.9752111086535323
$ws.Range("M15").Value = 13.32522833333333
$ws.Range("N15").Value = 39.975685
$ws.Range("O15").Value = 0.2709785829485105
$ws.Range("P15").Value = 0.3199146015909443
$ws.Range("Q15").Value = 23.23645321629667
$ws.Range("R15").Value = 139.41871929778
$ws.Range("S15").Value = 0.2643709786827099
$ws.Range("T15").Value = 0.3119842732919579

# Row 16
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 1.743794
$ws.Range("H16").Value = 3.487588
$ws.Range("I16").Value = 0.9756157693574768
$ws.Range("J16").Value = 0.9752111086535323
$ws.Range("M16").Value = 6.89049
$ws.Range("N16").Value = 20.67147
$ws.Range("O16").Value = 0.140123318663899
$ws.Range("P16").Value = 0.1654281868928364
$ws.Range("Q16").Value = 12.01559511906
$ws.Range("R16").Value = 72.09357071436
$ws.Range("S16").Value = 0.1367065193432027
$ws.Range("T16").Value = 0.1613274055423068

# Row 17
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 1.743794
$ws.Range("H17").Value = 3.487588
$ws.Range("I17").Value = 0.9756157693574768
$ws.Range("J17").Value = 0.9752111086535323
$ws.Range("M17").Value = 18.2696115
$ws.Range("N17").Value = 36.539223
$ws.Range("O17").Value = 0.3715263492262718
$ws.Range("P17").Value = 0.292413525083752
$ws.Range("Q17").Value = 31.858438916031
$ws.Range("R17").Value = 127.433755664124
$ws.Range("S17").Value = 0.3624669650369638
$ws.Range("T17").Value = 0.2851649179822133

# Row 18
$ws.Range("E18").Value = 2.0
$ws.Range("F18").Value = 1.0
$ws.Range("G18").Value = 1.743794
$ws.Range("H18").Value = 3.487588
$ws.Range("I18").Value = 0.9756157693574768
$ws.Range("J18").Value = 0.9752111086535323
$ws.Range("M18").Value = 1.355562
$ws.Range("N18").Value = 4.066686000000001
$ws.Range("O18").Value = 0.02756637715092428
$ws.Range("P18").Value = 0.03254458882907125
$ws.Range("Q18").Value = 2.363820882228
$ws.Range("R18").Value = 14.182925293368
$ws.Range("S18").Value = 0.02689419225249736
$ws.Range("T18").Value = 0.03173784455267194

# Row 19
$ws.Range("E19").Value = 2.0
$ws.Range("F19").Value = 1.0
$ws.Range("G19").Value = 1.743794
$ws.Range("H19").Value = 3.487588
$ws.Range("I19").Value = 0.9756157693574768
$ws.Range("J19").Value = 0.9752111086535323
$ws.Range("M19").Value = 5.037141666666667
$ws.Range("N19").Value = 15.111425
$ws.Range("O19").Value = 0.1024340804374633
$ws.Range("P19").Value = 0.1209326496430627
$ws.Range("Q19").Value = 8.783737415483333
$ws.Range("R19").Value = 52.7024244929
$ws.Range("S19").Value = 0.09993630419442144
$ws.Range("T19").Value = 0.1179348633308203
